$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Rating (C) and Year (D) columns so that
# numeric-looking values stay stored as text, matching the source data.
$ws.Range("C2:D21").NumberFormat = "@"

$ws.Range("A2").Value = 'The Fall Guy'
$ws.Range("C2").Value = '7.2'

$ws.Range("A3").Value = 'Mother of the Bride'
$ws.Range("B3").Value = 'Action'
$ws.Range("C3").Value = '4.8'

$ws.Range("A4").Value = 'Boy Kills World'
$ws.Range("B4").Value = 'Action'
$ws.Range("C4").Value = '6.6'

$ws.Range("A5").Value = 'Mothers'' Instinct'
$ws.Range("C5").Value = '6.3'

$ws.Range("A6").Value = 'Force of Nature: The Dry 2'
$ws.Range("B6").Value = 'Crime'
$ws.Range("C6").Value = '6'
$ws.Range("D6").Value = '2024'

$ws.Range("A7").Value = '[ES] The Courier'
$ws.Range("B7").Value = 'Crime'
$ws.Range("C7").Value = '6'
$ws.Range("D7").Value = '2024'

$ws.Range("A8").Value = 'Child''s Play'
$ws.Range("C8").Value = '6.7'
$ws.Range("D8").Value = '1988'

$ws.Range("A9").Value = 'The Last Stop in Yuma County'
$ws.Range("C9").Value = '7'
$ws.Range("D9").Value = '2023'

$ws.Range("A10").Value = 'Please Don''t Destroy: The Treasure of Foggy Mountain'
$ws.Range("C10").Value = '6.1'
$ws.Range("D10").Value = '2023'

$ws.Range("A11").Value = '[ML] Manjummel Boys'
$ws.Range("C11").Value = '8.4'
$ws.Range("D11").Value = '2024'

$ws.Range("A12").Value = 'Dune: Part Two'
$ws.Range("C12").Value = '8.6'
$ws.Range("D12").Value = '2024'

$ws.Range("A13").Value = 'Turtles All the Way Down'
$ws.Range("C13").Value = '6.7'
$ws.Range("D13").Value = '2024'

$ws.Range("A14").Value = 'Flynn'
$ws.Range("C14").Value = '4.5'
$ws.Range("D14").Value = '2024'

$ws.Range("A15").Value = 'The Idea of You'
$ws.Range("B15").Value = 'Action'
$ws.Range("C15").Value = '6.4'

$ws.Range("A16").Value = 'Unfrosted'

$ws.Range("A17").Value = 'One Bad Apple: A Hannah Swensen Mystery'
$ws.Range("C17").Value = '6.5'
$ws.Range("D17").Value = '2024'

$ws.Range("A18").Value = '[TR] About Dry Grasses'
$ws.Range("C18").Value = '7.8'

$ws.Range("A19").Value = 'The Hill'
$ws.Range("C19").Value = '6.7'
$ws.Range("D19").Value = '2023'

$ws.Range("A20").Value = 'Out of Darkness'
$ws.Range("C20").Value = '5.5'
$ws.Range("D20").Value = '2022'

$ws.Range("A21").Value = '[LT] Slow'
$ws.Range("C21").Value = '7.2'

